$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New employee names to insert (in nip order as shown in target sheet)
$data = @(
    @(6,  "Haryo Suro Kuncoro"),
    @(7,  "Muhammad Al-fatih Ritonga"),
    @(8,  "M. Irpandi"),
    @(9,  "Aditya Maulana"),
    @(11, "Ocha Sugiarto"),
    @(12, "Hafiizh Yaafi")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = 100
    $row++
}

$ws.Range("D6").Select()
